$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the values currently stored in A10 and G10 before we touch them
# (note: reading the COM `Value` property requires calling it, since the
# bare property reference resolves to an overload descriptor here).
$a10 = $ws.Range("A10").Value()
$g10 = $ws.Range("G10").Value()

# Clear the old placeholder content of row 4 and the old row 10 entirely.
$ws.Range("A4:H4").Clear()
$ws.Range("A10:H10").Clear()

# Write the relocated values into row 4: A4 gets what used to be in A10,
# H4 gets what used to be in G10 (now aligned under the H1 header).
$ws.Range("A4").Value = $a10
$ws.Range("B4:G4").Value = ""
$ws.Range("H4").Value = $g10

$ws.Range("A4:G4").Select()
